$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 -> year 2022, "Solar" column (E): 5 -> 6
$ws.Range("E24").Value = 6

# Row 26 -> year 2024, "Solar" column (E): 8 -> 14
$ws.Range("E26").Value = 14
